$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "43.697.18"
$ws.Range("E2").Value = "  -0.08%  "

Set-TextValue $ws.Range("D3") "2.327.05"
$ws.Range("E3").Value = "  +4.49%  "

$ws.Range("E4").Value = "  +0.11%  "

Set-TextValue $ws.Range("D5") "95.48"
$ws.Range("E5").Value = "  +6.83%  "

Set-TextValue $ws.Range("D6") "270.48"
$ws.Range("E6").Value = "  -0.22%  "

Set-TextValue $ws.Range("D7") "0.626"
$ws.Range("E7").Value = "  +0.95%  "

$ws.Range("E8").Value = "  +0.07%  "

Set-TextValue $ws.Range("D9") "0.621"
$ws.Range("E9").Value = "  +2.79%  "

Set-TextValue $ws.Range("D10") "44.94"
$ws.Range("E10").Value = "  -2.08%  "

$ws.Range("E11").Value = "  +2.97%  "

Set-TextValue $ws.Range("D12") "8.07"
$ws.Range("E12").Value = "  +3.83%  "

$ws.Range("E13").Value = "  +0.31%  "

Set-TextValue $ws.Range("D14") "2.680.27"

Set-TextValue $ws.Range("D15") "15.70"
$ws.Range("E15").Value = "  +4.56%  "

Set-TextValue $ws.Range("D16") "0.855"
$ws.Range("E16").Value = "  +8.26%  "

Set-TextValue $ws.Range("D17") "2.332.44"
$ws.Range("E17").Value = "  +4.80%  "

Set-TextValue $ws.Range("D18") "43.663.82"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("E19").Value = "  +6.00%  "

Set-TextValue $ws.Range("D20") "6.36"
$ws.Range("E20").Value = "  +6.66%  "

Set-TextValue $ws.Range("D21") "71.81"
$ws.Range("E21").Value = "  +2.12%  "

Set-TextValue $ws.Range("D22") "241.10"
$ws.Range("E22").Value = "  +3.93%  "

$ws.Range("E23").Value = "  -3.71%  "

$ws.Range("E24").Value = "  +9.85%  "

Set-TextValue $ws.Range("D25") "0.999"
$ws.Range("E25").Value = "  -0.07%  "

Set-TextValue $ws.Range("D26") "2.54"
$ws.Range("E26").Value = "  +1.21%  "

Set-TextValue $ws.Range("D27") "11.43"
$ws.Range("E27").Value = "  +4.76%  "

Set-TextValue $ws.Range("D28") "3.46"
$ws.Range("E28").Value = "  -2.88%  "

Set-TextValue $ws.Range("D29") "2.28"
$ws.Range("E29").Value = "  +0.64%  "

Set-TextValue $ws.Range("D30") "22.51"
$ws.Range("E30").Value = "  +8.74%  "

Set-TextValue $ws.Range("D31") "38.15"
$ws.Range("E31").Value = "  -1.00%  "

Set-TextValue $ws.Range("D32") "172.85"
$ws.Range("E32").Value = "  +0.18%  "

Set-TextValue $ws.Range("D33") "0.0899"
$ws.Range("E33").Value = "  -0.86%  "

$ws.Range("E34").Value = "  +3.10%  "

$ws.Range("E35").Value = "  +2.55%  "

Set-TextValue $ws.Range("D36") "0.0358"
$ws.Range("E36").Value = "  +2.20%  "

$ws.Range("E37").Value = "  -2.97%  "

Set-TextValue $ws.Range("D38") "4.38"
$ws.Range("E38").Value = "  +3.14%  "

Set-TextValue $ws.Range("D39") "3.35"
$ws.Range("E39").Value = "  -2.54%  "

Set-TextValue $ws.Range("D40") "2.34"
$ws.Range("E40").Value = "  +8.81%  "

$ws.Range("E41").Value = "  +10.21%  "

Set-TextValue $ws.Range("D42") "1.36"
$ws.Range("E42").Value = "  +18.37%  "

Set-TextValue $ws.Range("D43") "12.11"
$ws.Range("E43").Value = "  -1.80%  "

Set-TextValue $ws.Range("D44") "9.15"
$ws.Range("E44").Value = "  +7.89%  "

Set-TextValue $ws.Range("D45") "61.83"
$ws.Range("E45").Value = "  -2.10%  "

$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("E47").Value = "  +4.36%  "

Set-TextValue $ws.Range("D48") "100.28"
$ws.Range("E48").Value = "  +0.38%  "

Set-TextValue $ws.Range("D49") "1.22"
$ws.Range("E49").Value = "  +3.20%  "

Set-TextValue $ws.Range("D50") "2.560.87"
$ws.Range("E50").Value = "  +4.59%  "

Set-TextValue $ws.Range("D51") "0.181"
$ws.Range("E51").Value = "  +13.20%  "
